# Scheduled runner update: refresh market-board derived price/profit figures
# across the crafting-job leve sheets (currentAveragePrice* / LevePrice* /
# LeveProfit* columns) with latest pulled data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 960.3415
$ws.Range("J17").Value = 960.3415
$ws.Range("L17").Value = 2881.0245
$ws.Range("N17").Value = -3217.0245

$ws.Range("H40").Value = 4245.609
$ws.Range("J40").Value = 4577.778
$ws.Range("L40").Value = 4577.778
$ws.Range("N40").Value = -4927.778

$ws.Range("H64").Value = 4964.647
$ws.Range("J64").Value = 4981.1875
$ws.Range("L64").Value = 4981.1875
$ws.Range("N64").Value = -5477.1875

$ws.Range("H67").Value = 4964.647
$ws.Range("J67").Value = 4981.1875
$ws.Range("L67").Value = 4981.1875
$ws.Range("N67").Value = -6697.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 8750
$ws.Range("I33").Value = 7500
$ws.Range("J33").Value = 10000
$ws.Range("K33").Value = 7500
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = -7164
$ws.Range("N33").Value = -10672

$ws.Range("H92").Value = 68694.25
$ws.Range("J92").Value = 68694.25
$ws.Range("L92").Value = 68694.25
$ws.Range("N92").Value = -73686.25

$ws.Range("H94").Value = 2025.625
$ws.Range("I94").Value = 2368
$ws.Range("K94").Value = 2368
$ws.Range("M94").Value = -1917

$ws.Range("H134").Value = 46001.914
$ws.Range("I134").Value = 701.5
$ws.Range("K134").Value = 2104.5
$ws.Range("M134").Value = 430.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H58").Value = 1843.7333
$ws.Range("I58").Value = 1801.25
$ws.Range("K58").Value = 1801.25
$ws.Range("M58").Value = -1598.25

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H136").Value = 1843.7333
$ws.Range("I136").Value = 1801.25
$ws.Range("K136").Value = 5403.75
$ws.Range("M136").Value = -2853.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 398.33334
$ws.Range("I6").Value = 278
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 834
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -721
$ws.Range("N6").Value = -3226

$ws.Range("H7").Value = 148
$ws.Range("I7").Value = 160.25
$ws.Range("J7").Value = 99
$ws.Range("K7").Value = 480.75
$ws.Range("L7").Value = 297
$ws.Range("M7").Value = -368.75
$ws.Range("N7").Value = -521

$ws.Range("H112").Value = 10587.883
$ws.Range("I112").Value = 9642.5
$ws.Range("J112").Value = 14999.667
$ws.Range("K112").Value = 28927.5
$ws.Range("L112").Value = 44999.001
$ws.Range("M112").Value = -27819.5
$ws.Range("N112").Value = -47215.001

$ws.Range("H115").Value = 23006.2
$ws.Range("J115").Value = 39015.5
$ws.Range("L115").Value = 117046.5
$ws.Range("N115").Value = -119396.5

$ws.Range("H131").Value = 8938.454
$ws.Range("I131").Value = 17299.666
$ws.Range("J131").Value = 5803
$ws.Range("K131").Value = 51898.99800000001
$ws.Range("L131").Value = 17409
$ws.Range("M131").Value = -46858.99800000001
$ws.Range("N131").Value = -27489

$ws.Range("H134").Value = 4660.1113
$ws.Range("I134").Value = 3992.1333
$ws.Range("J134").Value = 8000
$ws.Range("K134").Value = 11976.3999
$ws.Range("L134").Value = 24000
$ws.Range("M134").Value = -6906.3999
$ws.Range("N134").Value = -34140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H64").Value = 30120.334
$ws.Range("J64").Value = 30120.334
$ws.Range("L64").Value = 30120.334
$ws.Range("N64").Value = -30616.334

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H67").Value = 30120.334
$ws.Range("J67").Value = 30120.334
$ws.Range("L67").Value = 30120.334
$ws.Range("N67").Value = -31836.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 12506.5
$ws.Range("I32").Value = 12506.5
$ws.Range("K32").Value = 12506.5
$ws.Range("M32").Value = -12189.5

$ws.Range("H46").Value = 5510.9585
$ws.Range("I46").Value = 3357.125
$ws.Range("J46").Value = 9818.625
$ws.Range("K46").Value = 3357.125
$ws.Range("L46").Value = 9818.625
$ws.Range("M46").Value = -3169.125
$ws.Range("N46").Value = -10194.625

$ws.Range("H54").Value = 39999
$ws.Range("J54").Value = 39999
$ws.Range("L54").Value = 39999
$ws.Range("N54").Value = -41287

$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("M61").Value = -798

$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("M113").Value = 1170

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 23500
$ws.Range("J59").Value = 23500
$ws.Range("L59").Value = 23500
$ws.Range("N59").Value = -24976

$ws.Range("H75").Value = 25356000
$ws.Range("J75").Value = 38003000
$ws.Range("L75").Value = 38003000
$ws.Range("N75").Value = -38004872

$ws.Range("H78").Value = 25356000
$ws.Range("J78").Value = 38003000
$ws.Range("L78").Value = 114009000
$ws.Range("N78").Value = -114018360

$ws.Range("H132").Value = 1502.9117
$ws.Range("I132").Value = 1342.6786
$ws.Range("K132").Value = 4028.0358
$ws.Range("M132").Value = -1498.0358
